# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" (quarterly fund-position detail) right
#   before the existing "2021-Q2" sheet, populated with the new quarter's
#   fund holdings.
# - Update the "总计" (totals) summary sheet: the row that used to describe
#   2021-Q2 now describes 2022-Q3 (with its own numbers), and a new row is
#   appended at the bottom preserving the original 2021-Q2 totals.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$oldQuarter = $wb.Worksheets.Item("2021-Q2")

# ---------------------------------------------------------------------------
# 1. Duplicate the "2021-Q2" sheet (same layout/styling as the sheet we are
#    about to add) right before itself, then rename the duplicate -- this is
#    the cleanest way to get a new quarter tab that keeps the sheet-level
#    look & feel (outline settings, page margins, column widths, ...) other
#    quarter tabs already have.
# ---------------------------------------------------------------------------
$oldQuarter.Copy($oldQuarter)
$newQuarter = $wb.Worksheets.Item("2021-Q2 (2)")
$newQuarter.Name = "2022-Q3"

# Page setup on the duplicated sheet was inherited from "2021-Q2"; line it up
# with the rest of the workbook (same margins as the "总计" sheet) instead.
$newQuarter.PageSetup.LeftMargin = 54
$newQuarter.PageSetup.RightMargin = 54
$newQuarter.PageSetup.TopMargin = 72
$newQuarter.PageSetup.BottomMargin = 72
$newQuarter.PageSetup.HeaderMargin = 36
$newQuarter.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2. Overwrite the duplicated sheet with the 2022-Q3 fund-position table.
# ---------------------------------------------------------------------------
$newQuarter.Range("B1").Value = "基金代码"
$newQuarter.Range("C1").Value = "基金名称"
$newQuarter.Range("D1").Value = "基金规模"
$newQuarter.Range("E1").Value = "股票总仓位"
$newQuarter.Range("F1").Value = "仓位占比"
$newQuarter.Range("G1").Value = "持有市值(亿元)"
$newQuarter.Range("H1").Value = "仓位排名"

$newQuarter.Range("A2").Value = 0
$newQuarter.Range("B2").Value = "'005585"
$newQuarter.Range("C2").Value = "银河文体娱乐主题灵活配置混合A"
$newQuarter.Range("D2").Value = "'3.01"
$newQuarter.Range("E2").Value = "'90.28"
$newQuarter.Range("F2").Value = "'4.87"
$newQuarter.Range("G2").Value = "'0.1466"
$newQuarter.Range("H2").Value = 6

$newQuarter.Range("A3").Value = 1
$newQuarter.Range("B3").Value = "'015667"
$newQuarter.Range("C3").Value = "银河文体娱乐主题灵活配置混合C"
$newQuarter.Range("D3").Value = "'0.41"
$newQuarter.Range("E3").Value = "'90.28"
$newQuarter.Range("F3").Value = "'4.87"
$newQuarter.Range("G3").Value = "'0.0200"
$newQuarter.Range("H3").Value = 6

# The apostrophe-forced text entries above pick up a quote-prefix style;
# clean that back off by re-pasting the plain (unstyled) format used
# throughout these tables.
$summary.Range("C2").Copy()
$newQuarter.Range("B2:G3").PasteSpecial(-4122)

# Re-apply the bold/bordered look (matching the "总计" sheet) to the header
# row and the leading index column.
$summary.Range("B1:D1").Copy()
$newQuarter.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$newQuarter.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Update the "总计" sheet: row 2 now reflects 2022-Q3, and a new row 3 is
#    appended for the original 2021-Q2 totals.
# ---------------------------------------------------------------------------
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("D2").Value = 0.17

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.01

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Leave "2021-Q2" as the selected/active tab, same as before the edit
#    (re-fetch by name since the earlier Copy() shuffled sheet positions).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
